# "fixed export and fixing maps"
#
# The sheet held a small area table for წყალტუბო municipality with an
# (unused/obsolete) census-reference note under the title, and a
# three-year comparison (1989 / 2002 / 2014) of the area figure. The
# fix drops the stale census note line and collapses the table down to
# just the current (2014) figure, and gives the worksheet its proper
# name instead of the generic "1".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "1" to "წყალტუბო".
$ws.Name = "წყალტუბო"

# Drop the 1989 and 2002 columns - only the 2014 figure is kept (it
# slides left into column B).
$ws.Range("B:C").Delete()

# Remove the now-unneeded "(census results)" note row right under the
# title; everything below shifts up one row.
$ws.Range("2:2").Delete()

# Leave the selection where the author left it.
$ws.Range("A2").Select()
